$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01575462090081432
$ws.Range("C2").Value = 0.1965624703810414

$ws.Range("B3").Value = 0.05458100723811542
$ws.Range("C3").Value = 0.2361503079336973

$ws.Range("B4").Value = 0.941907158318113
$ws.Range("C4").Value = 0.4912571846679999

$ws.Range("B5").Value = 0.9956299508590023
$ws.Range("C5").Value = 0.4235547453023541

$ws.Range("B6").Value = 0.9923357471225299
$ws.Range("C6").Value = 0.8144106704366469

$ws.Range("B7").Value = 0.9850802293349668
$ws.Range("C7").Value = 0.3398143927567758

$ws.Range("B8").Value = 0.002559327483177185
$ws.Range("C8").Value = 0.1780865287780762
